$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so values like "1.010" / "2.086.53"
# are stored as literal text (matching the inlineStr cells in the source) rather
# than being auto-parsed as numbers by Excel's smart input.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "29.619.33"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "2.086.53"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "342.51"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "0.5155"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").Value = "0.4389"
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("D9").Value = "0.09224"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").Value = "51.98"
$ws.Range("E10").Value = "  -3.53%  "
$ws.Range("D11").Value = "1.175"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "25.05"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "2.091.69"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "6.732"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "8.150"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "99.79"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "0.00001155"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "21.06"
$ws.Range("E19").Value = "  +8.83%  "
$ws.Range("D22").Value = "6.169"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").Value = "29.672.28"
$ws.Range("E23").Value = "  -3.45%  "
$ws.Range("D24").Value = "12.62"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D26").Value = "2.335.30"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "21.84"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "162.95"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").Value = "2.513"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "132.36"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").Value = "1.135"
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("D32").Value = "0.1049"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").Value = "1.626"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "6.174"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").Value = "3.963"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").Value = "6.020"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").Value = "10.31"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").Value = "0.02567"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "0.06701"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "12.43"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("D41").Value = "0.2230"
$ws.Range("E41").Value = "  -3.45%  "
$ws.Range("D42").Value = "0.6815"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "1.289"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "0.6605"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("D45").Value = "14.20"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("D46").Value = "2.308"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").Value = "3.608"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("D48").Value = "1.215"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "0.00000000339"
$ws.Range("E49").Value = "  -5.89%  "
$ws.Range("D50").Value = "81.40"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "1.167"
$ws.Range("E51").Value = "  -2.09%  "

$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E25").Value = "  -3.02%  "

# Restore the original (default) cell style now that the text values are set,
# so we don't leave a lingering custom number format on the cells.
$colD.Style = "Normal"
